$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update some existing quantities in the "capacitores nano" block ---
$ws.Range("B48").Value = 104
$ws.Range("B50").Value = 103
$ws.Range("B51").Value = 102
$ws.Range("C51").Value = 4
$ws.Range("C53").Value = 4

# --- Insert a new row "12p" / 1 right after "15p" (row 53) ---
$ws.Rows("54:54").Insert()
$ws.Range("B54").Value = "12p"
$ws.Range("C54").Value = 1

# Update "22p" quantity (now on row 56 after the insert above)
$ws.Range("C56").Value = 8

# --- Insert a new row "0,1u" / 3 right after "1u" (now row 62) ---
$ws.Rows("63:63").Insert()
$ws.Range("B63").Value = "0,1u"
$ws.Range("C63").Value = 3

# Update "0,33u" quantity (now on row 66 after both inserts)
$ws.Range("C66").Value = 2

# --- Column A width (auto-sized after editing, widened slightly) ---
$ws.Columns("A:A").ColumnWidth = 15.3

# --- Refresh the view: scroll down and move the selection ---
$excel.ActiveWindow.ScrollRow = 45
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C67").Select()
